$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 5022
$ws.Range("I34").Value = 5022
$ws.Range("K34").Value = 5022
$ws.Range("M34").Value = -4819
$ws.Range("H36").Value = 5022
$ws.Range("I36").Value = 5022
$ws.Range("K36").Value = 5022
$ws.Range("M36").Value = -4307
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -6850
$ws.Range("H43").Value = 2317.3333
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2317.3333
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2317.3333
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2455.3333
$ws.Range("H64").Value = 2714.2856
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 2916.6667
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 2916.6667
$ws.Range("M64").Value = -1252
$ws.Range("N64").Value = -3412.6667
$ws.Range("H67").Value = 2714.2856
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 2916.6667
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 2916.6667
$ws.Range("M67").Value = -642
$ws.Range("N67").Value = -4632.6667
$ws.Range("H80").Value = 1085.625
$ws.Range("I80").Value = 858.75
$ws.Range("J80").Value = 1312.5
$ws.Range("K80").Value = 2576.25
$ws.Range("L80").Value = 3937.5
$ws.Range("M80").Value = -1578.25
$ws.Range("N80").Value = -5933.5
$ws.Range("H83").Value = 1085.625
$ws.Range("I83").Value = 858.75
$ws.Range("J83").Value = 1312.5
$ws.Range("K83").Value = 7728.75
$ws.Range("L83").Value = 11812.5
$ws.Range("M83").Value = -2736.75
$ws.Range("N83").Value = -21796.5
$ws.Range("H86").Value = 23166.666
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 23166.666
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws.Range("H92").Value = 909.2353000000001
$ws.Range("I92").Value = 702.25
$ws.Range("J92").Value = 1406
$ws.Range("K92").Value = 702.25
$ws.Range("L92").Value = 1406
$ws.Range("M92").Value = 545.75
$ws.Range("N92").Value = -3902
$ws.Range("H100").Value = 715.8889
$ws.Range("J100").Value = 1000
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082
$ws.Range("H113").Value = 2473.875
$ws.Range("I113").Value = 2298.5
$ws.Range("K113").Value = 2298.5
$ws.Range("M113").Value = 955.5
$ws.Range("H141").Value = 599
$ws.Range("I141").Value = 599
$ws.Range("K141").Value = 1797
$ws.Range("M141").Value = 3383

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 885.6667
$ws.Range("I32").Value = 835.75
$ws.Range("K32").Value = 835.75
$ws.Range("M32").Value = -548.75
$ws.Range("H33").Value = 13552.5
$ws.Range("H74").Value = 2299.8
$ws.Range("I74").Value = 1899.75
$ws.Range("K74").Value = 1899.75
$ws.Range("M74").Value = -1025.75
$ws.Range("H77").Value = 2299.8
$ws.Range("I77").Value = 1899.75
$ws.Range("K77").Value = 9498.75
$ws.Range("M77").Value = -5130.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 29500
$ws.Range("J76").Value = 29500
$ws.Range("L76").Value = 29500
$ws.Range("N76").Value = -30130
$ws.Range("H79").Value = 29500
$ws.Range("J79").Value = 29500
$ws.Range("L79").Value = 29500
$ws.Range("N79").Value = -31684
$ws.Range("H99").Value = 4349.778
$ws.Range("I99").Value = 4349.778
$ws.Range("K99").Value = 4349.778
$ws.Range("M99").Value = -2851.778

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 13.666667
$ws.Range("H41").Value = 26029.5
$ws.Range("I41").Value = 7059
$ws.Range("J41").Value = 45000
$ws.Range("K41").Value = 7059
$ws.Range("L41").Value = 45000
$ws.Range("M41").Value = -6631
$ws.Range("N41").Value = -45856
$ws.Range("H58").Value = 1556.8
$ws.Range("I58").Value = 1556.8
$ws.Range("K58").Value = 1556.8
$ws.Range("M58").Value = -1353.8
$ws.Range("H59").Value = 11001
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H132").Value = 1847.4
$ws.Range("I132").Value = 1846.75
$ws.Range("K132").Value = 5540.25
$ws.Range("M132").Value = -3010.25
$ws.Range("H136").Value = 1556.8
$ws.Range("I136").Value = 1556.8
$ws.Range("K136").Value = 4670.4
$ws.Range("M136").Value = -2120.4

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 5730
$ws.Range("I18").Value = 5730
$ws.Range("K18").Value = 17190
$ws.Range("M18").Value = -17021
$ws.Range("H23").Value = 129.83333
$ws.Range("J23").Value = 129.83333
$ws.Range("L23").Value = 389.49999
$ws.Range("N23").Value = -859.49999
$ws.Range("H111").Value = 350
$ws.Range("I111").Value = 350
$ws.Range("K111").Value = 1050
$ws.Range("M111").Value = 2017
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H131").Value = 2515
$ws.Range("I131").Value = 30
$ws.Range("K131").Value = 90
$ws.Range("M131").Value = 4950
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 7497
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17697

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5100000
$ws.Range("J3").Value = 6666666.5
$ws.Range("L3").Value = 6666666.5
$ws.Range("N3").Value = -6666898.5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H97").Value = 736.25
$ws.Range("I97").Value = 648.3333
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 648.3333
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -152.3333
$ws.Range("N97").Value = -1992
$ws.Range("H123").Value = 75799.60000000001
$ws.Range("J123").Value = 75799.60000000001
$ws.Range("L123").Value = 75799.60000000001
$ws.Range("N123").Value = -80699.60000000001
$ws.Range("H132").Value = 1312.5
$ws.Range("I132").Value = 1312.5
$ws.Range("K132").Value = 3937.5
$ws.Range("M132").Value = -1407.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4737.5
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 5983.3335
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 5983.3335
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -6359.3335
$ws.Range("H55").Value = 343
$ws.Range("I55").Value = 367.6
$ws.Range("J55").Value = 312.25
$ws.Range("K55").Value = 367.6
$ws.Range("L55").Value = 312.25
$ws.Range("M55").Value = -194.6
$ws.Range("N55").Value = -658.25
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 4003.5
$ws.Range("I122").Value = 4003.5
$ws.Range("K122").Value = 12010.5
$ws.Range("M122").Value = -9560.5
$ws.Range("H132").Value = 3024.75
$ws.Range("I132").Value = 3004
$ws.Range("J132").Value = 3031.6667
$ws.Range("K132").Value = 9012
$ws.Range("L132").Value = 9095.000100000001
$ws.Range("M132").Value = -6482
$ws.Range("N132").Value = -14155.0001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 10273
$ws.Range("I82").Value = 10273
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 10273
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -9890
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 10273
$ws.Range("I85").Value = 10273
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 10273
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -8947
$ws.Range("N85").ClearContents()
$ws.Range("H119").Value = 41250
$ws.Range("J119").Value = 41250
$ws.Range("L119").Value = 41250
$ws.Range("N119").Value = -50926
$ws.Range("H131").Value = 56000
$ws.Range("J131").Value = 56000
$ws.Range("L131").Value = 56000
$ws.Range("N131").Value = -66080
$ws.Range("H132").Value = 2179.5557
$ws.Range("I132").Value = 2160.7646
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 6482.293799999999
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -3952.293799999999
$ws.Range("N132").Value = -12557
$ws.Range("H136").Value = 2524.75
$ws.Range("I136").Value = 2524.75
$ws.Range("K136").Value = 7574.25
$ws.Range("M136").Value = -5024.25
